# "10.02/2024 - return to host 34"
# Updates Date_of_sales (column I) from 45349 to 45350 for rows 2-30,
# and reorders words in the Param text columns (C and G) so that
# "б/к" / "H" move from the front/middle to the end of the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column I: Date_of_sales, rows 2-30: 45349 -> 45350 ---
for ($r = 2; $r -le 30; $r++) {
    $ws.Range("I$r").Value = 45350
}

# --- Column C: reorder descriptive words ---
$ws.Range("C3").Value  = "сер легк б/к"
$ws.Range("C4").Value  = "сер легк б/к"
$ws.Range("C10").Value = "210B H сер C Type"
$ws.Range("C12").Value = "202B сер C Type"
$ws.Range("C13").Value = "202B LS-2 H сер C Type"
$ws.Range("C14").Value = "сер груз б/к"
$ws.Range("C15").Value = "сер легк б/к"
$ws.Range("C16").Value = "сер легк б/к"

# --- Column G: reorder descriptive words (comma separated) ---
$ws.Range("G4").Value  = "сер, легк, б/к"
$ws.Range("G5").Value  = "сер, легк, б/к"
$ws.Range("G11").Value = "210B, H, сер, C, Type"
$ws.Range("G12").Value = "210B, H, сер, C, Type"
$ws.Range("G14").Value = "202B, сер, C, Type"
$ws.Range("G15").Value = "202B, LS-2, H, сер, C, Type"
$ws.Range("G16").Value = "202B, LS-2, H, сер, C, Type"
$ws.Range("G17").Value = "202B, LS-2, H, сер, C, Type"
$ws.Range("G18").Value = "сер, груз, б/к"
$ws.Range("G19").Value = "сер, груз, б/к"
$ws.Range("G20").Value = "сер, груз, б/к"
$ws.Range("G21").Value = "сер, груз, б/к"
$ws.Range("G22").Value = "сер, легк, б/к"
$ws.Range("G23").Value = "сер, легк, б/к"
